$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round Q2 and R2 to nearest whole number
$ws.Range("Q2").Value = 772263
$ws.Range("R2").Value = 7120316

# Clear Starttid (Z2) and Sluttid (AB2) cells entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
